$wb = $excel.ActiveWorkbook

# Rename the "temperature_c" sheet to "temperature"
$ws = $wb.Worksheets.Item("temperature_c")
$ws.Name = "temperature"

# Make it the active/selected sheet (moves tabSelected + activeTab onto it,
# and away from whatever was active before - "genotype")
$ws.Activate()
